# Auto-generated edit script: Add 2024-03-29 violent-crime data update
# Updates 2024 (column K) year-to-date totals (and a few upstream J/B corrections)
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$sheetEdits = @(
    @{ Sheet = 'Citywide Totals'; Cells = @(@{Cell="K2"; Value=1694}, @{Cell="J3"; Value=8077}, @{Cell="K3"; Value=1606}, @{Cell="B4"; Value=1695}, @{Cell="J4"; Value=1799}, @{Cell="K4"; Value=350}, @{Cell="K5"; Value=103}, @{Cell="K6"; Value=2101}, @{Cell="B7"; Value=23328}, @{Cell="K7"; Value=5854}) }
    @{ Sheet = 'Austin'; Cells = @(@{Cell="K2"; Value=111}, @{Cell="K3"; Value=109}, @{Cell="K6"; Value=126}, @{Cell="K7"; Value=373}) }
    @{ Sheet = 'South Chicago'; Cells = @(@{Cell="K2"; Value=51}, @{Cell="K3"; Value=40}, @{Cell="K7"; Value=120}) }
    @{ Sheet = 'Garfield Park'; Cells = @(@{Cell="K4"; Value=16}, @{Cell="K7"; Value=240}) }
    @{ Sheet = 'West Pullman'; Cells = @(@{Cell="K6"; Value=28}, @{Cell="K7"; Value=100}) }
    @{ Sheet = 'Grand Crossing'; Cells = @(@{Cell="K2"; Value=45}, @{Cell="K5"; Value=11}, @{Cell="K6"; Value=71}, @{Cell="K7"; Value=199}) }
    @{ Sheet = 'New City'; Cells = @(@{Cell="K6"; Value=66}, @{Cell="K7"; Value=148}) }
    @{ Sheet = 'Woodlawn'; Cells = @(@{Cell="K2"; Value=35}, @{Cell="K4"; Value=4}, @{Cell="K6"; Value=34}, @{Cell="K7"; Value=111}) }
    @{ Sheet = 'By Neighborhood'; Cells = @(@{Cell="K4"; Value=25}, @{Cell="K6"; Value=49}, @{Cell="K7"; Value=164}, @{Cell="K8"; Value=373}, @{Cell="B11"; Value=410}, @{Cell="K14"; Value=34}, @{Cell="K15"; Value=53}, @{Cell="K19"; Value=159}, @{Cell="K25"; Value=29}, @{Cell="K29"; Value=271}, @{Cell="K31"; Value=66}, @{Cell="K33"; Value=240}, @{Cell="K36"; Value=64}, @{Cell="K37"; Value=199}, @{Cell="K42"; Value=206}, @{Cell="J48"; Value=323}, @{Cell="K48"; Value=64}, @{Cell="K49"; Value=35}, @{Cell="K50"; Value=28}, @{Cell="K52"; Value=154}, @{Cell="K54"; Value=101}, @{Cell="B63"; Value=400}, @{Cell="K63"; Value=22}, @{Cell="K64"; Value=39}, @{Cell="K65"; Value=148}, @{Cell="K73"; Value=57}, @{Cell="K76"; Value=86}, @{Cell="K79"; Value=156}, @{Cell="K80"; Value=20}, @{Cell="K83"; Value=120}, @{Cell="J85"; Value=1194}, @{Cell="K85"; Value=297}, @{Cell="K88"; Value=72}, @{Cell="K91"; Value=53}, @{Cell="K95"; Value=100}, @{Cell="K99"; Value=111}, @{Cell="B101"; Value=23328}, @{Cell="K101"; Value=5854}) }
    @{ Sheet = 'Gage Park'; Cells = @(@{Cell="K5"; Value=2}, @{Cell="K7"; Value=66}) }
    @{ Sheet = 'Lincoln Park'; Cells = @(@{Cell="K6"; Value=24}, @{Cell="K7"; Value=35}) }
    @{ Sheet = 'Loop'; Cells = @(@{Cell="K4"; Value=6}, @{Cell="K6"; Value=40}, @{Cell="K7"; Value=101}) }
    @{ Sheet = 'Englewood'; Cells = @(@{Cell="K2"; Value=74}, @{Cell="K3"; Value=88}, @{Cell="K4"; Value=13}, @{Cell="K6"; Value=90}, @{Cell="K7"; Value=271}) }
    @{ Sheet = 'Lake View'; Cells = @(@{Cell="J3"; Value=62}, @{Cell="K4"; Value=10}, @{Cell="K6"; Value=29}, @{Cell="J7"; Value=323}, @{Cell="K7"; Value=64}) }
    @{ Sheet = 'Chatham'; Cells = @(@{Cell="K2"; Value=49}, @{Cell="K7"; Value=159}) }
    @{ Sheet = 'River North'; Cells = @(@{Cell="K6"; Value=49}, @{Cell="K7"; Value=86}) }
    @{ Sheet = 'Bridgeport'; Cells = @(@{Cell="K6"; Value=13}, @{Cell="K7"; Value=34}) }
    @{ Sheet = 'Ashburn'; Cells = @(@{Cell="K6"; Value=17}, @{Cell="K7"; Value=49}) }
    @{ Sheet = 'Humboldt Park'; Cells = @(@{Cell="K2"; Value=50}, @{Cell="K3"; Value=55}, @{Cell="K6"; Value=90}, @{Cell="K7"; Value=206}) }
    @{ Sheet = 'Washington Park'; Cells = @(@{Cell="K3"; Value=20}, @{Cell="K7"; Value=53}) }
    @{ Sheet = 'Roseland'; Cells = @(@{Cell="K3"; Value=55}, @{Cell="K7"; Value=156}) }
    @{ Sheet = 'Near South Side'; Cells = @(@{Cell="K3"; Value=13}, @{Cell="K7"; Value=39}) }
    @{ Sheet = 'Grand Boulevard'; Cells = @(@{Cell="K2"; Value=28}, @{Cell="K7"; Value=64}) }
    @{ Sheet = 'Auburn Gresham'; Cells = @(@{Cell="K6"; Value=40}, @{Cell="K7"; Value=164}) }
    @{ Sheet = 'East Side'; Cells = @(@{Cell="K2"; Value=12}, @{Cell="K7"; Value=29}) }
    @{ Sheet = 'Brighton Park'; Cells = @(@{Cell="K2"; Value=16}, @{Cell="K6"; Value=22}, @{Cell="K7"; Value=53}) }
    @{ Sheet = 'Lincoln Square'; Cells = @(@{Cell="K6"; Value=18}, @{Cell="K7"; Value=28}) }
    @{ Sheet = 'Belmont Cragin'; Cells = @(@{Cell="B4"; Value=29}, @{Cell="B7"; Value=410}) }
    @{ Sheet = 'Portage Park'; Cells = @(@{Cell="K6"; Value=25}, @{Cell="K7"; Value=57}) }
    @{ Sheet = 'United Center'; Cells = @(@{Cell="K6"; Value=41}, @{Cell="K7"; Value=72}) }
    @{ Sheet = 'South Shore'; Cells = @(@{Cell="K2"; Value=109}, @{Cell="K3"; Value=95}, @{Cell="J4"; Value=73}, @{Cell="K6"; Value=73}, @{Cell="J7"; Value=1194}, @{Cell="K7"; Value=297}) }
    @{ Sheet = 'Rush & Division'; Cells = @(@{Cell="K6"; Value=10}, @{Cell="K7"; Value=20}) }
    @{ Sheet = 'Little Village'; Cells = @(@{Cell="K3"; Value=33}, @{Cell="K7"; Value=154}) }
    @{ Sheet = 'Archer Heights'; Cells = @(@{Cell="K2"; Value=8}, @{Cell="K7"; Value=25}) }
)

foreach ($entry in $sheetEdits) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($c in $entry.Cells) {
        $ws.Range($c.Cell).Value = $c.Value
    }
}
